$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3186.2727
$ws.Range("I74").Value = 3099.7837
$ws.Range("J74").Value = 3643.4285
$ws.Range("K74").Value = 3099.7837
$ws.Range("L74").Value = 3643.4285
$ws.Range("M74").Value = -2163.7837
$ws.Range("N74").Value = -5515.4285
# Row 77
$ws.Range("H77").Value = 3186.2727
$ws.Range("I77").Value = 3099.7837
$ws.Range("J77").Value = 3643.4285
$ws.Range("K77").Value = 15498.9185
$ws.Range("L77").Value = 18217.1425
$ws.Range("M77").Value = -10818.9185
$ws.Range("N77").Value = -27577.1425
# Row 82
$ws.Range("H82").Value = 4034.9
$ws.Range("I82").Value = 1724.8334
$ws.Range("K82").Value = 5174.5002
$ws.Range("M82").Value = -4768.5002
# Row 85
$ws.Range("H85").Value = 4034.9
$ws.Range("I85").Value = 1724.8334
$ws.Range("K85").Value = 5174.5002
$ws.Range("M85").Value = -3770.5002
# Row 132
$ws.Range("H132").Value = 5849558.5
$ws.Range("I132").Value = 2029.6154
$ws.Range("J132").Value = 18519204
$ws.Range("K132").Value = 6088.8462
$ws.Range("L132").Value = 55557612
$ws.Range("M132").Value = -3558.8462
$ws.Range("N132").Value = -55562672

$ws = $wb.Worksheets.Item("ARM")
# Row 86
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32372
# Row 89
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 90000
$ws.Range("N89").Value = -101856
# Row 98
$ws.Range("H98").Value = 34000
$ws.Range("J98").Value = 34000
$ws.Range("L98").Value = 34000
$ws.Range("N98").Value = -39990

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1956.6703
$ws.Range("I86").Value = 1967.816
$ws.Range("J86").Value = 1714.25
$ws.Range("K86").Value = 1967.816
$ws.Range("L86").Value = 1714.25
$ws.Range("M86").Value = -844.816
$ws.Range("N86").Value = -3960.25
# Row 89
$ws.Range("H89").Value = 1956.6703
$ws.Range("I89").Value = 1967.816
$ws.Range("J89").Value = 1714.25
$ws.Range("K89").Value = 9839.08
$ws.Range("L89").Value = 8571.25
$ws.Range("M89").Value = -4223.08
$ws.Range("N89").Value = -19803.25

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 9348.947
$ws.Range("I86").Value = 14803.444
$ws.Range("J86").Value = 4439.9
$ws.Range("K86").Value = 14803.444
$ws.Range("L86").Value = 4439.9
$ws.Range("M86").Value = -13680.444
$ws.Range("N86").Value = -6685.9
# Row 89
$ws.Range("H89").Value = 9348.947
$ws.Range("I89").Value = 14803.444
$ws.Range("J89").Value = 4439.9
$ws.Range("K89").Value = 74017.22
$ws.Range("L89").Value = 22199.5
$ws.Range("M89").Value = -68401.22
$ws.Range("N89").Value = -33431.5

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1474.3684
$ws.Range("I34").Value = 630
$ws.Range("J34").Value = 1966.9166
$ws.Range("K34").Value = 1890
$ws.Range("L34").Value = 5900.7498
$ws.Range("M34").Value = -1806
$ws.Range("N34").Value = -6068.7498
# Row 92
$ws.Range("H92").Value = 2439416
$ws.Range("I92").Value = 510.4
$ws.Range("J92").Value = 4878321.5
$ws.Range("K92").Value = 1531.2
$ws.Range("L92").Value = 14634964.5
$ws.Range("M92").Value = -283.1999999999998
$ws.Range("N92").Value = -14637460.5

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 16787.555
$ws.Range("I13").Value = 147
$ws.Range("J13").Value = 30100
$ws.Range("K13").Value = 147
$ws.Range("L13").Value = 30100
$ws.Range("M13").Value = -8
$ws.Range("N13").Value = -30378
# Row 70
$ws.Range("H70").Value = 3122561.8
$ws.Range("I70").Value = 1921676
$ws.Range("J70").Value = 5057322
$ws.Range("K70").Value = 1921676
$ws.Range("L70").Value = 5057322
$ws.Range("M70").Value = -1921406
$ws.Range("N70").Value = -5057862
# Row 73
$ws.Range("H73").Value = 3122561.8
$ws.Range("I73").Value = 1921676
$ws.Range("J73").Value = 5057322
$ws.Range("K73").Value = 1921676
$ws.Range("L73").Value = 5057322
$ws.Range("M73").Value = -1920740
$ws.Range("N73").Value = -5059194
# Row 92
$ws.Range("H92").Value = 3252.6
$ws.Range("J92").Value = 3252.6
$ws.Range("L92").Value = 3252.6
$ws.Range("N92").Value = -6996.6
# Row 107
$ws.Range("H107").Value = 405.25
$ws.Range("J107").Value = 450
$ws.Range("L107").Value = 450
$ws.Range("N107").Value = -4290
# Row 140
$ws.Range("H140").Value = 43926.145
$ws.Range("J140").Value = 43926.145
$ws.Range("L140").Value = 43926.145
$ws.Range("N140").Value = -54286.145
# Row 141
$ws.Range("H141").Value = 46060
$ws.Range("J141").Value = 46060
$ws.Range("L141").Value = 46060
$ws.Range("N141").Value = -56420

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 573.5
$ws.Range("I46").Value = 746.2
$ws.Range("J46").Value = 400.8
$ws.Range("K46").Value = 746.2
$ws.Range("L46").Value = 400.8
$ws.Range("M46").Value = -558.2
$ws.Range("N46").Value = -776.8
# Row 82
$ws.Range("H82").Value = 1789.3077
$ws.Range("I82").Value = 1820
$ws.Range("K82").Value = 1820
$ws.Range("M82").Value = -1459
# Row 85
$ws.Range("H85").Value = 1789.3077
$ws.Range("I85").Value = 1820
$ws.Range("K85").Value = 1820
$ws.Range("M85").Value = -572
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 93
$ws.Range("H93").Value = 18594.268
$ws.Range("J93").Value = 22328.273
$ws.Range("L93").Value = 22328.273
$ws.Range("N93").Value = -24824.273
# Row 94
$ws.Range("H94").Value = 31665
$ws.Range("J94").Value = 31665
$ws.Range("L94").Value = 31665
$ws.Range("N94").Value = -33017
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 97
$ws.Range("H97").Value = 16769.416
$ws.Range("J97").Value = 16769.416
$ws.Range("L97").Value = 16769.416
$ws.Range("N97").Value = -18751.416
# Row 98
$ws.Range("H98").Value = 23500
$ws.Range("J98").Value = 23500
$ws.Range("L98").Value = 23500
$ws.Range("N98").Value = -29490
# Row 99
$ws.Range("H99").Value = 26000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 26000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 26000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -31990
# Row 100
$ws.Range("H100").Value = 1341.2122
$ws.Range("I100").Value = 1080.5555
$ws.Range("J100").Value = 1654
$ws.Range("K100").Value = 1080.5555
$ws.Range("L100").Value = 1654
$ws.Range("M100").Value = -539.5554999999999
$ws.Range("N100").Value = -2736
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 103
$ws.Range("H103").Value = 20932.834
$ws.Range("J103").Value = 20932.834
$ws.Range("L103").Value = 20932.834
$ws.Range("N103").Value = -23276.834
# Row 104
$ws.Range("H104").Value = 23566.666
$ws.Range("J104").Value = 23566.666
$ws.Range("L104").Value = 23566.666
$ws.Range("N104").Value = -30554.666
# Row 105
$ws.Range("H105").Value = 35312.5
$ws.Range("J105").Value = 35312.5
$ws.Range("L105").Value = 35312.5
$ws.Range("N105").Value = -42300.5
# Row 106
$ws.Range("H106").Value = 34184.5
$ws.Range("J106").Value = 34184.5
$ws.Range("L106").Value = 34184.5
$ws.Range("N106").Value = -36708.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 15909.523
$ws.Range("I81").Value = 2125.25
$ws.Range("J81").Value = 19152.883
$ws.Range("K81").Value = 4250.5
$ws.Range("L81").Value = 38305.766
$ws.Range("M81").Value = -3189.5
$ws.Range("N81").Value = -40427.766
# Row 84
$ws.Range("H84").Value = 15909.523
$ws.Range("I84").Value = 2125.25
$ws.Range("J84").Value = 19152.883
$ws.Range("K84").Value = 21252.5
$ws.Range("L84").Value = 191528.83
$ws.Range("M84").Value = -15948.5
$ws.Range("N84").Value = -202136.83
# Row 135
$ws.Range("H135").Value = 44686
$ws.Range("J135").Value = 44686
$ws.Range("L135").Value = 44686
$ws.Range("N135").Value = -54826
# Row 140
$ws.Range("H140").Value = 38400
$ws.Range("J140").Value = 38400
$ws.Range("L140").Value = 38400
$ws.Range("N140").Value = -48760
# Row 141
$ws.Range("H141").Value = 52696
$ws.Range("J141").Value = 52696
$ws.Range("L141").Value = 52696
$ws.Range("N141").Value = -63056
